# Task_任务.xlsx — stage-reward rebalance + precondition cleanup
#
# Context (from the commit "feat：update stage config change and excel"):
#   - World-1 stage rewards ("rewards" column G) are bumped:
#       Stage_stageName_1 (简单/easy)   id 1001  1|50  -> 1|100
#       Stage_stageName_2 (普通/normal) id 1002  1|75  -> 1|150
#       Stage_stageName_3 (困难/hard)   id 1003  1|100 -> 1|200
#   - The "preconditions" column F is cleared out for the stage rows
#     (ids 1002..1015 / rows 6-19) — these stages no longer declare an
#     explicit precondition id.
#   - The active worksheet selection moves from I64 to I12.
#
# (All the other cell-index churn visible in the raw OOXML diff — the
#  B/C/D/G shared-string index bumps on rows 52-82 — is just fallout of
#  inserting one new shared string ("1|50") into the table; the actual
#  cell *values* there are unchanged, so Excel's own save path reproduces
#  that automatically and nothing further needs to be written for them.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- World 1 stage rewards bump -------------------------------------------------
$ws.Range("G5").Value = "1|100"   # id 1001, Stage_stageName_1, 通关世界1·简单
$ws.Range("G6").Value = "1|150"   # id 1002, Stage_stageName_2, 通关世界1·普通
$ws.Range("G7").Value = "1|200"   # id 1003, Stage_stageName_3, 通关世界1·困难

# --- Clear stale "preconditions" values on the stage rows (ids 1002-1015) ------
$ws.Range("F6:F19").ClearContents()

# --- Move the saved worksheet selection (cosmetic, matches the author's last
#     edit position) -------------------------------------------------------------
$ws.Range("I12").Select() | Out-Null
